# Reorder the comma-separated "Recorded By" list in column G:
#   - entries WITHOUT an "@" (e.g. "System", "system") move to the front,
#     keeping their original relative order
#   - entries WITH an "@" (actual email addresses) move to the back,
#     in reversed relative order
# This matches the transformation observed between the original and the
# edited workbook for every row in the "Recorded By" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Reorder-RecordedBy($s) {
    $nonAtItems = @()
    $atItems = @()

    $parts = $s.Split(",")
    foreach ($p in $parts) {
        $t = $p.Trim()
        if ($t.Length -eq 0) {
            continue
        }
        if ($t.Contains("@")) {
            $atItems += $t
        } else {
            $nonAtItems += $t
        }
    }

    if ($atItems.Count -gt 1) {
        $atItems = $atItems[($atItems.Count - 1)..0]
    }

    $result = $nonAtItems + $atItems
    return ($result -join ", ")
}

for ($row = 2; $row -le 157; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $newVal = Reorder-RecordedBy $val
        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
